$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the lecture date cells (column B) to the new schedule
$ws.Range("B2").Value = "27.08 *08:15 - 12:00* (**Aud M**)"
$ws.Range("B3").Value = "03.09 *08:15 - 12:00* (**Aud M**)"
$ws.Range("B4").Value = "10.09 *08:15 - 12:00* (**Aud M**)"
$ws.Range("B5").Value = "17.09 *08:15 - 12:00* (**Aud M**)"

# Seminar dates (column F) are cleared out to "-"
$ws.Range("F2").Value = "-"
$ws.Range("F3").Value = "-"
$ws.Range("F4").Value = "-"
$ws.Range("F5").Value = "-"
$ws.Range("F6").Value = "-"

# Assignment due dates (column G) are cleared out to "-"
$ws.Range("G2").Value = "-"
$ws.Range("G3").Value = "-"
$ws.Range("G4").Value = "-"
$ws.Range("G5").Value = "-"

# Move the active selection to B6 (matches the saved selection state)
$ws.Range("B6").Select() | Out-Null
